$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: score_psi_test_segments
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("score_psi_test_segments")

$ws.Range("D6").Value = 24.05913978494624
$ws.Range("F6").Value = 179
$ws.Range("G6").Value = -1.754891769266587
$ws.Range("H6").Value = 0.0012355007655722
$ws.Range("J6").Value = 0.2653365456990043
$ws.Range("K6").Value = 0.005423885254925055

$ws.Range("D7").Value = 4.56989247311828
$ws.Range("F7").Value = 34
$ws.Range("G7").Value = 0.8102415835580246
$ws.Range("H7").Value = 0.00158129640982946
$ws.Range("I7").Value = 0.3204213058735216
$ws.Range("J7").Value = 0.315759486682328
$ws.Range("K7").Value = 0.007005181664754515

# ---------------------------------------------------------------------------
# Sheet: score_psi_oot_segments
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("score_psi_oot_segments")

$ws.Range("J2").Value = 0.1416258266621568

$ws.Range("D3").Value = 20.6734534064213
$ws.Range("F3").Value = 264
$ws.Range("G3").Value = -1.683042061856646
$ws.Range("H3").Value = 0.001317257622436057
$ws.Range("J3").Value = 0.1520087400391462
$ws.Range("K3").Value = 0.001430496954768251

$ws.Range("D4").Value = 24.66718872357087
$ws.Range("F4").Value = 315
$ws.Range("G4").Value = 0.9679607947356894
$ws.Range("H4").Value = 0.0003874892063125338
$ws.Range("J4").Value = 0.1842496862321489
$ws.Range("K4").Value = 0.001817986161080785

$ws.Range("D5").Value = 3.680501174628034
$ws.Range("F5").Value = 47
$ws.Range("G5").Value = 0.3908066462628111
$ws.Range("H5").Value = 0.0004386969513771168
$ws.Range("J5").Value = 0.20306586887407
$ws.Range("K5").Value = 0.002256683112457901

$ws.Range("D6").Value = 26.54659357870008
$ws.Range("F6").Value = 339
$ws.Range("G6").Value = 0.7325620244872544
$ws.Range("H6").Value = 0.0002049945495771488
$ws.Range("J6").Value = 0.266404408955256
$ws.Range("K6").Value = 0.00246167766203505

$ws.Range("I7").Value = 0.3204213058735216
$ws.Range("J7").Value = 0.315173418310558
$ws.Range("K7").Value = 0.002477477534896982

# ---------------------------------------------------------------------------
# Sheet: calibration_risk_bands_segments
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("calibration_risk_bands_segments")

$ws.Range("F3").Value = 26.36360920517067

# ---------------------------------------------------------------------------
# Sheet: calibration_deciles_segments
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("calibration_deciles_segments")

$ws.Range("G2").Value = 8.719101405670601
$ws.Range("K2").Value = 0.007211321588540948

$ws.Range("E5").Value = 15.26187809239385
$ws.Range("G5").Value = 8.360169151700639
$ws.Range("K5").Value = 0.01299428876467782
